$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 299 (shifts existing rows 299:370 down to 300:371)
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the new weekly data point
$ws.Range("A299").Value = 4
$ws.Range("B299").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C299").Value = "Los Lagos"
$ws.Range("D299").Value = 44943
$ws.Range("E299").Value = 10
$ws.Range("F299").Value = 100112043
$ws.Range("G299").Value = "Pepino ensalada"
$ws.Range("H299").Value = "Sin especificar"
$ws.Range("I299").Value = "Primera"
$ws.Range("J299").Value = 400
$ws.Range("K299").Value = 16000
$ws.Range("L299").Value = 18000
$ws.Range("M299").Value = 17000
$ws.Range("N299").Value = "$/caja 60 unidades"
$ws.Range("O299").Value = "Región de Arica y Parinacota"
$ws.Range("P299").Value = 283
$ws.Range("Q299").Value = 60
$ws.Range("R299").Value = "Hortaliza"
